# Updates the cryptos list prices (column D) and Volume(1h) percentages (column E) in Sheet1
# to reflect the latest GitHub Actions refresh of cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 2; D = '39.882.72'; DNumericLike = $false; E = '  +1.17%  ' }
    @{ Row = 3; D = '2.221.62'; DNumericLike = $false; E = '  +0.08%  ' }
    @{ Row = 4; D = $null; DNumericLike = $false; E = '  +0.07%  ' }
    @{ Row = 5; D = '292.22'; DNumericLike = $true; E = '  -1.56%  ' }
    @{ Row = 6; D = '87.50'; DNumericLike = $true; E = '  +6.83%  ' }
    @{ Row = 8; D = $null; DNumericLike = $false; E = '  -0.04%  ' }
    @{ Row = 9; D = '0.473'; DNumericLike = $true; E = '  +0.50%  ' }
    @{ Row = 10; D = '30.33'; DNumericLike = $true; E = '  +1.69%  ' }
    @{ Row = 11; D = $null; DNumericLike = $false; E = '  +1.40%  ' }
    @{ Row = 12; D = '47.52'; DNumericLike = $true; E = '  +1.54%  ' }
    @{ Row = 13; D = $null; DNumericLike = $false; E = '  +1.78%  ' }
    @{ Row = 14; D = '6.41'; DNumericLike = $true; E = '  +2.09%  ' }
    @{ Row = 15; D = '2.562.99'; DNumericLike = $false; E = '  +0.27%  ' }
    @{ Row = 16; D = '14.09'; DNumericLike = $true; E = '  +0.16%  ' }
    @{ Row = 17; D = '2.220.24'; DNumericLike = $false; E = '  +0.01%  ' }
    @{ Row = 18; D = '0.731'; DNumericLike = $true; E = '  +2.05%  ' }
    @{ Row = 19; D = '39.830.04'; DNumericLike = $false; E = '  +1.21%  ' }
    @{ Row = 20; D = '11.58'; DNumericLike = $true; E = '  +12.13%  ' }
    @{ Row = 21; D = '0.0₃0884'; DNumericLike = $false; E = '  +1.12%  ' }
    @{ Row = 22; D = '5.84'; DNumericLike = $true; E = '  +1.45%  ' }
    @{ Row = 23; D = '65.81'; DNumericLike = $true; E = '  +1.45%  ' }
    @{ Row = 24; D = '235.92'; DNumericLike = $true; E = '  +2.36%  ' }
    @{ Row = 25; D = $null; DNumericLike = $false; E = '  +0.08%  ' }
    @{ Row = 26; D = '2.48'; DNumericLike = $true; E = '  +2.87%  ' }
    @{ Row = 27; D = $null; DNumericLike = $false; E = '  +1.33%  ' }
    @{ Row = 28; D = '22.84'; DNumericLike = $true; E = '  +0.57%  ' }
    @{ Row = 29; D = $null; DNumericLike = $false; E = '  +1.27%  ' }
    @{ Row = 30; D = '9.27'; DNumericLike = $true; E = '  +1.61%  ' }
    @{ Row = 31; D = '32.82'; DNumericLike = $true; E = '  +2.65%  ' }
    @{ Row = 32; D = '152.33'; DNumericLike = $true; E = '  +2.09%  ' }
    @{ Row = 33; D = $null; DNumericLike = $false; E = '  -0.03%  ' }
    @{ Row = 34; D = '4.96'; DNumericLike = $true; E = '  +2.94%  ' }
    @{ Row = 35; D = '0.0722'; DNumericLike = $true; E = '  +3.22%  ' }
    @{ Row = 36; D = $null; DNumericLike = $false; E = '  +1.34%  ' }
    @{ Row = 37; D = '2.82'; DNumericLike = $true; E = '  +6.47%  ' }
    @{ Row = 38; D = $null; DNumericLike = $false; E = '  +1.53%  ' }
    @{ Row = 39; D = '15.98'; DNumericLike = $true; E = '  +1.55%  ' }
    @{ Row = 40; D = '0.0991'; DNumericLike = $true; E = '  +2.92%  ' }
    @{ Row = 41; D = '1.71'; DNumericLike = $true; E = '  +2.57%  ' }
    @{ Row = 42; D = '2.096.45'; DNumericLike = $false; E = '  +9.55%  ' }
    @{ Row = 43; D = '3.81'; DNumericLike = $true; E = '  +4.46%  ' }
    @{ Row = 44; D = $null; DNumericLike = $false; E = '  +5.89%  ' }
    @{ Row = 45; D = $null; DNumericLike = $false; E = '  +3.62%  ' }
    @{ Row = 46; D = $null; DNumericLike = $false; E = '  +9.24%  ' }
    @{ Row = 47; D = '17.64'; DNumericLike = $true; E = '  +6.93%  ' }
    @{ Row = 48; D = '2.64'; DNumericLike = $true; E = '  +0.72%  ' }
    @{ Row = 49; D = '2.433.80'; DNumericLike = $false; E = '  +0.28%  ' }
    @{ Row = 50; D = '70.96'; DNumericLike = $true; E = '  -0.79%  ' }
    @{ Row = 51; D = '89.40'; DNumericLike = $true; E = '  +1.45%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)
        if ($u.DNumericLike) {
            # Value looks like a plain number to Excel's parser (e.g. "292.22").
            # The source data is text, so force a quote-prefixed literal via
            # Formula instead of Value, keeping the cell a text value like the original.
            $cell.Formula = "'" + $u.D
        } else {
            # Not parseable as a number (e.g. "39.882.72" or "0.0₃0884") -
            # Value assignment already stores it as text with no style changes.
            $cell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
